$d = $word.ActiveDocument

# 1. Update activation date
$d.Content.Find.Execute(
    "Ativação: 01/01/2020",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ativação: 01/01/2025",
    2
)

# 2. Append to Portuguese short program summary
$d.Content.Find.Execute(
    "Determinação das suscetibilidades e vocações do meio ambiente e o conceito de sustentabilidade ambiental.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Determinação das suscetibilidades e vocações do meio ambiente e o conceito de sustentabilidade ambiental; Gestão de Recursos Hídricos.",
    2
)

# 3. Append to English short program summary
$d.Content.Find.Execute(
    "Environment susceptibilities and vocations determination and environmental susceptibility concept.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Environment susceptibilities and vocations determination and environmental susceptibility concept, Water Resources Management.",
    2
)

# 4. Append to Portuguese full program
$d.Content.Find.Execute(
    "Estrutura institucional e marcos legais em recursos hídricos no Brasil.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Estrutura institucional e marcos legais em recursos hídricos no Brasil. Legislação e instrumentos pertinentes. A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina.",
    2
)

# 5. Append to English full program
$d.Content.Find.Execute(
    "Brazilian institutional structure and legal frameworks in water resources.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Brazilian institutional structure and legal frameworks in water resources; Legislation and relevant instrument. The discipline may have didactic trips to complement the content of the discipline.",
    2
)
